# Update cryptos list values (price & 1h volume change) as scraped on Fri May  3 05:34:05 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.685.82"
$ws.Range("E2").Value = "  +3.96%  "
$ws.Range("D3").Value = "3.006.13"
$ws.Range("E3").Value = "  +3.39%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'564.18"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").Value = "'140.62"
$ws.Range("E6").Value = "  +9.16%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "2.995.49"
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("E10").Value = "  +7.48%  "
$ws.Range("D11").Value = "'5.29"
$ws.Range("E11").Value = "  +12.64%  "
$ws.Range("D12").Value = "'0.453"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("E13").Value = "  +7.65%  "
$ws.Range("D14").Value = "'33.99"
$ws.Range("E14").Value = "  +4.88%  "
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "3.503.80"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("E17").Value = "  +5.08%  "
$ws.Range("D18").Value = "3.001.65"
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("D19").Value = "59.659.26"
$ws.Range("E19").Value = "  +3.87%  "
$ws.Range("D20").Value = "'439.30"
$ws.Range("E20").Value = "  +6.20%  "
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("D22").Value = "'0.719"
$ws.Range("E22").Value = "  +5.26%  "
$ws.Range("D23").Value = "'13.52"
$ws.Range("E23").Value = "  +3.71%  "
$ws.Range("D24").Value = "'7.10"
$ws.Range("E24").Value = "  +2.49%  "
$ws.Range("D25").Value = "'80.55"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'2.23"
$ws.Range("E27").Value = "  +12.49%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +4.33%  "
$ws.Range("E30").Value = "  +6.52%  "
$ws.Range("E31").Value = "  +6.47%  "
$ws.Range("E32").Value = "  +3.28%  "
$ws.Range("D33").Value = "'0.105"
$ws.Range("E33").Value = "  +9.60%  "
$ws.Range("D34").Value = "0.0₃0787"
$ws.Range("E34").Value = "  +16.82%  "
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  +8.06%  "
$ws.Range("D36").Value = "'5.92"
$ws.Range("E36").Value = "  +5.31%  "
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("D38").Value = "'49.10"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'2.84"
$ws.Range("E39").Value = "  +13.25%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "'8.59"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'401.29"
$ws.Range("E41").Value = "  +7.87%  "
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").Value = "2.767.22"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("D44").Value = "'0.106"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +7.37%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'123.44"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "'34.44"
$ws.Range("E48").Value = "  +22.68%  "
$ws.Range("E49").Value = "  +5.63%  "
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("E51").Value = "  +4.04%  "
